# Apply the coupon-sheet update:
#  - rename the existing coupon's title from "كوبون نمشي" to "نمشي"
#  - append a new coupon row (نوون / scaios2026) with the same
#    description, link (+hyperlink), image, countries and note as row 2

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 2: update the title in A2 ---
$ws.Range("A2").Value = "نمشي"

# --- Row 3: new coupon entry ---
$ws.Range("A3").Value = "نوون"
$ws.Range("B3").Value = "خصم على جميع المنتجات"
$ws.Range("C3").Value = "scaios2026"
$ws.Range("D3").Value = "https://www.discountcoupon.online"
$ws.Range("E3").Value = "https://f.top4top.io/p_3389y71vl1.png"
$ws.Range("F3").Value = "السعودية"
$ws.Range("G3").Value = "صالح للاستخدام عدة مرات"

# Add the hyperlink on D3 (mirrors the one on D2) and match its style
$ws.Hyperlinks.Add($ws.Range("D3"), "https://www.discountcoupon.online/")
$ws.Range("D3").Style = $ws.Range("D2").Style

# Update the last active selection to match the saved workbook state
$ws.Range("G12").Select()
